# "2024-3" performance-summary edits, sheet "绩效表" (first sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 15: business-type correction ---
$ws.Cells.Item(15, 4).Value = "固定业务"

# --- Row 17: new entry for 黄礼闯 (#11 - 戴心怡沟通) ---
$ws.Cells.Item(17, 1).Value = "黄礼闯"
$ws.Cells.Item(17, 2).Value = 11
$ws.Cells.Item(17, 3).Value = "戴心怡沟通"
$ws.Cells.Item(17, 4).Value = "其他业务"
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 9).Value = "完成"
$ws.Cells.Item(17, 11).Value = 0.03
$ws.Cells.Item(17, 12).Value = 180
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0

# Columns E/G/H/J on this row are genuine empty-string text (matching rows
# 7/15/16 above), not merely blank cells. A bare `.Value = ""` clears a
# cell outright (real Excel semantics), so force text with a leading
# apostrophe and then restore the plain (non quote-prefixed) number format
# by pulling it from the equivalent cell one row up - that keeps the
# column's original style index (37/38/39/23) intact.
$ws.Cells.Item(17, 5).Value = "'"
$ws.Cells.Item(16, 5).Copy()
$ws.Cells.Item(17, 5).PasteSpecial(-4122)

$ws.Cells.Item(17, 7).Value = "'"
$ws.Cells.Item(16, 7).Copy()
$ws.Cells.Item(17, 7).PasteSpecial(-4122)

$ws.Cells.Item(17, 8).Value = "'"
$ws.Cells.Item(16, 8).Copy()
$ws.Cells.Item(17, 8).PasteSpecial(-4122)

$ws.Cells.Item(17, 10).Value = "'"
$ws.Cells.Item(16, 1).Copy()
$ws.Cells.Item(17, 10).PasteSpecial(-4122)

# --- Row 29: "基本业务" summary totals ---
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 8).Value = 0.226
$ws.Cells.Item(29, 9).Value = "0.226=0.226"
$ws.Cells.Item(29, 10).Value = 1.154

# --- Row 30: "其他业务" summary totals ---
$ws.Cells.Item(30, 8).Value = 0.928
$ws.Cells.Item(30, 9).Value = "0.254+0.008+0.058+0.058+0.17+0.114+0.086+0.142+0.008+0.03=0.928"
